$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 196.125
$ws.Range("I8").Value = 196.125
$ws.Range("K8").Value = 588.375
$ws.Range("M8").Value = -449.375
$ws.Range("H28").Value = 875.1053000000001
$ws.Range("I28").Value = 671.75
$ws.Range("K28").Value = 671.75
$ws.Range("M28").Value = -186.75
$ws.Range("H38").Value = 769987.3
$ws.Range("J38").Value = 1500
$ws.Range("L38").Value = 4500
$ws.Range("N38").Value = -5244
$ws.Range("H51").Value = 6437.25
$ws.Range("J51").Value = 15000
$ws.Range("L51").Value = 15000
$ws.Range("N51").Value = -15968
$ws.Range("H87").Value = 93450.664
$ws.Range("J87").Value = 93450.664
$ws.Range("L87").Value = 93450.664
$ws.Range("N87").Value = -95946.664
$ws.Range("H90").Value = 93450.664
$ws.Range("J90").Value = 93450.664
$ws.Range("L90").Value = 280351.992
$ws.Range("N90").Value = -292831.992
$ws.Range("H95").Value = 31499
$ws.Range("J95").Value = 31499
$ws.Range("L95").Value = 31499
$ws.Range("N95").Value = -36991
$ws.Range("H98").Value = 1557.3334
$ws.Range("I98").Value = 676
$ws.Range("K98").Value = 676
$ws.Range("M98").Value = 822
$ws.Range("H107").Value = 1143.8889
$ws.Range("I107").Value = 753.3333
$ws.Range("J107").Value = 1925
$ws.Range("K107").Value = 753.3333
$ws.Range("L107").Value = 1925
$ws.Range("M107").Value = 1166.6667
$ws.Range("N107").Value = -5765
$ws.Range("H122").Value = 1557.3334
$ws.Range("I122").Value = 676
$ws.Range("K122").Value = 2028
$ws.Range("M122").Value = 422

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2863.8125
$ws.Range("I32").Value = 1383.7142
$ws.Range("K32").Value = 1383.7142
$ws.Range("M32").Value = -1096.7142
$ws.Range("H45").Value = 1272.8
$ws.Range("I45").Value = 1341.125
$ws.Range("K45").Value = 1341.125
$ws.Range("M45").Value = -964.125
$ws.Range("H95").Value = 65539
$ws.Range("J95").Value = 65539
$ws.Range("L95").Value = 65539
$ws.Range("N95").Value = -71031

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 855.5
$ws.Range("I107").Value = 862.4286
$ws.Range("K107").Value = 862.4286
$ws.Range("M107").Value = 1057.5714

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("L42").ClearContents()
$ws.Range("M42").ClearContents()
$ws.Range("N42").Value = 0
$ws.Range("H56").Value = 58833.332
$ws.Range("J56").Value = 58833.332
$ws.Range("L56").Value = 58833.332
$ws.Range("N56").Value = -60523.332
$ws.Range("H62").Value = 2966.6667
$ws.Range("I62").Value = 2950
$ws.Range("J62").Value = 3000
$ws.Range("K62").Value = 2950
$ws.Range("L62").Value = 3000
$ws.Range("M62").Value = -2326
$ws.Range("N62").Value = -4248
$ws.Range("H65").Value = 2966.6667
$ws.Range("I65").Value = 2950
$ws.Range("J65").Value = 3000
$ws.Range("K65").Value = 14750
$ws.Range("L65").Value = 15000
$ws.Range("M65").Value = -11630
$ws.Range("N65").Value = -21240
$ws.Range("H134").Value = 2659.3845
$ws.Range("I134").Value = 2659.3845
$ws.Range("K134").Value = 7978.1535
$ws.Range("M134").Value = -5443.1535

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").ClearContents()
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = 0
$ws.Range("H132").Value = 874.75
$ws.Range("I132").Value = 874.75
$ws.Range("K132").Value = 7872.75
$ws.Range("M132").Value = -5342.75
$ws.Range("H134").Value = 10655.9
$ws.Range("I134").Value = 725.1111
$ws.Range("J134").Value = 100033
$ws.Range("K134").Value = 2175.3333
$ws.Range("L134").Value = 300099
$ws.Range("M134").Value = 2894.6667
$ws.Range("N134").Value = -310239
$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("L135").ClearContents()
$ws.Range("M135").ClearContents()
$ws.Range("N135").Value = 0

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 15394.8
$ws.Range("J92").Value = 15394.8
$ws.Range("L92").Value = 15394.8
$ws.Range("N92").Value = -19138.8

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4433.636
$ws.Range("I7").Value = 4433.636
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 4433.636
$ws.Range("L7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -4321.636
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("L19").ClearContents()
$ws.Range("M19").ClearContents()
$ws.Range("N19").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()
$ws.Range("H61").Value = 1800.2
$ws.Range("I61").Value = 1800.2
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1800.2
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -1598.2
$ws.Range("H113").Value = 1800.2
$ws.Range("I113").Value = 1800.2
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1800.2
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = 369.8
$ws.Range("H122").Value = 4167.2
$ws.Range("I122").Value = 4096.857
$ws.Range("K122").Value = 12290.571
$ws.Range("M122").Value = -9840.571
$ws.Range("H126").Value = 4433.636
$ws.Range("I126").Value = 4433.636
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 13300.908
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -10830.908
$ws.Range("H132").Value = 4674.5
$ws.Range("I132").Value = 4933
$ws.Range("K132").Value = 14799
$ws.Range("M132").Value = -12269

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 28428.5
$ws.Range("J46").Value = 28428.5
$ws.Range("L46").Value = 28428.5
$ws.Range("N46").Value = -28890.5
$ws.Range("H96").Value = 1432.6666
$ws.Range("I96").Value = 1449.5
$ws.Range("J96").Value = 1399
$ws.Range("K96").Value = 1449.5
$ws.Range("L96").Value = 1399
$ws.Range("M96").Value = -76.5
$ws.Range("N96").Value = -4145
$ws.Range("H122").Value = 3541.652
$ws.Range("I122").Value = 4067.125
$ws.Range("J122").Value = 2340.5715
$ws.Range("K122").Value = 12201.375
$ws.Range("L122").Value = 7021.7145
$ws.Range("M122").Value = -9751.375
$ws.Range("N122").Value = -11921.7145
$ws.Range("H126").Value = 3930.5
$ws.Range("I126").Value = 3608.3333
$ws.Range("J126").Value = 4897
$ws.Range("K126").Value = 10824.9999
$ws.Range("L126").Value = 14691
$ws.Range("M126").Value = -8354.999899999999
$ws.Range("N126").Value = -19631
$ws.Range("H132").Value = 4632.8335
$ws.Range("I132").Value = 4632.8335
$ws.Range("K132").Value = 13898.5005
$ws.Range("M132").Value = -11368.5005
$ws.Range("H134").Value = 28428.5
$ws.Range("J134").Value = 28428.5
$ws.Range("L134").Value = 85285.5
$ws.Range("N134").Value = -90355.5
$ws.Range("H135").Value = 190357.5
$ws.Range("J135").Value = 190357.5
$ws.Range("L135").Value = 190357.5
$ws.Range("N135").Value = -200497.5
